$d = $word.ActiveDocument

# Locate the paragraph that ends the "favicon" bullet, after which the
# new note about the whatsapp footer message should be inserted.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*en todas las secciones del sitio web.*") {
        $target = $para
    }
}

if ($target -eq $null) {
    throw "Could not locate the paragraph ending in 'en todas las secciones del sitio web.'"
}

# Create a new paragraph right after it; Word copies the paragraph mark's
# run formatting (Arial) automatically, same as the surrounding paragraphs.
$tail = $target.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

# The freshly inserted (still empty) paragraph is the next one.
$newPara = $target.Next()
$newRange = $newPara.Range
$newRange.Collapse(0)

# Build the curly-quoted fragment without backtick escapes (PowerShell
# treats `a` as the BEL control character, not the letter "a").
$ldquo = [char]0x201C
$rdquo = [char]0x201D
$quotedA = $ldquo + "a" + $rdquo
$middleText = " en el enlace de etiqueta " + $quotedA + " de "

$rFonts = '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
          '<w:p>' +
            '<w:pPr><w:rPr>' + $rFonts + '</w:rPr></w:pPr>' +
            '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t xml:space="preserve">Se personaliza mensaje de </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>whatspp</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t xml:space="preserve">' + $middleText + '</w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>whatsapp</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t xml:space="preserve"> en </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>footer</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>.</w:t></w:r>' +
          '</w:p>' +
        '</w:body>' +
      '</w:document>' +
    '</pkg:xmlData>' +
  '</pkg:part>' +
'</pkg:package>'

$newRange.InsertXML($xml)

Write-Output "Inserted paragraph about whatsapp footer message."
